$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C8").Value = "Indica que pretende editar lista de Materiais"

$ws.Range("C9").Select()
